$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2 through 267). The recorded value 45172 (2023-09-03) is being
# refreshed to 45175 (2023-09-06) for all of them.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 267 }

$ws.Range("C2:C$lastRow").Value = 45175
